$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.269.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.833.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4972"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.42%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1005"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +27.62%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.16"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.459"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.71"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.000"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.825.21"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.339"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001147"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.13"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06653"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.035"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.305.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.35"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.235"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.26"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.83"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.039.72"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.439"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.25"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.049"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.615"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.599"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06790"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.47%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02357"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2155"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.48"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.996"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6239"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.178"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.20"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5952"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.688"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.29"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.953"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.121"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.49%  "
